# Attendance sheet: add three date columns (B, C, D) with PRESENT/ABSENT
# status for each person, matching the "Database and error solved" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): date labels in B1:D1, styled like A1 ("Name").
# The values look like dates (2/7/23, 4/4/23, 4/3/23) but must be kept as
# literal text, so they are entered with a leading apostrophe (the normal
# Excel technique for forcing text entry instead of date auto-conversion).
$ws.Range("B1").Value = "'2/7/23"
$ws.Range("C1").Value = "'4/4/23"
$ws.Range("D1").Value = "'4/3/23"

$ws.Range("B1:D1").Font.Bold = $true
$ws.Range("B1:D1").Borders.LineStyle = 1
$ws.Range("B1:D1").HorizontalAlignment = -4108
$ws.Range("B1:D1").VerticalAlignment = -4160

# Attendance status per person (rows 2-8) across the three dates, in row order.
$statusByRow = @(
    @("PRESENT", "PRESENT", "PRESENT"),  # row 2: Sumit
    @("PRESENT", "PRESENT", "PRESENT"),  # row 3: Nigel
    @("PRESENT", "PRESENT", "PRESENT"),  # row 4: Abhay
    @("PRESENT", "PRESENT", "PRESENT"),  # row 5: Afzal
    @("PRESENT", "PRESENT", "PRESENT"),  # row 6: Prakhar
    @("ABSENT",  "ABSENT",  "ABSENT"),   # row 7: Avon
    @("ABSENT",  "ABSENT",  "ABSENT")    # row 8: Kainaat
)

for ($i = 0; $i -lt $statusByRow.Length; $i++) {
    $row = $i + 2
    $vals = $statusByRow[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}
